$wb = $excel.ActiveWorkbook

# sheet1 (index 1)
$ws = $wb.Worksheets.Item(1)
$ws.Range("F4").Value = 145
$ws.Range("F5").Value = 2042
$ws.Range("G5").Value = "不可售"
$ws.Range("F6").Value = 4173
$ws.Range("F7").Value = 542
$ws.Range("F8").Value = 1047
$ws.Range("F9").Value = 659
$ws.Range("F10").Value = 373
$ws.Range("F11").Value = 100
$ws.Range("F12").Value = 2195
$ws.Range("F13").Value = 398
$ws.Range("F14").Value = 657923
$ws.Range("F15").Value = 1624
$ws.Range("F16").Value = 505
$ws.Range("F17").Value = 1452
$ws.Range("F18").Value = 670
$ws.Range("F19").Value = 541
$ws.Range("F20").Value = 1269
$ws.Range("F21").Value = 2222
$ws.Range("F22").Value = 1134
$ws.Range("F23").Value = 2691
$ws.Range("F24").Value = 1547
$ws.Range("F25").Value = 800
$ws.Range("F26").Value = 1532
$ws.Range("F27").Value = 25
$ws.Range("F28").Value = 527
$ws.Range("F29").Value = 1083
$ws.Range("F30").Value = 272
$ws.Range("F31").Value = 1078
$ws.Range("F32").Value = 42
$ws.Range("F34").Value = 2012
$ws.Range("F35").Value = 1367
$ws.Range("F36").Value = 573
$ws.Range("F37").Value = 1222
$ws.Range("F38").Value = 2439
$ws.Range("F39").Value = 1142
$ws.Range("F40").Value = 25
$ws.Range("F41").Value = 197
$ws.Range("F42").Value = 2575
$ws.Range("F43").Value = 207
$ws.Range("F44").Value = 980
$ws.Range("F45").Value = 3114
$ws.Range("F46").Value = 1006
$ws.Range("F48").Value = 877
$ws.Range("F49").Value = 150

# sheet2 (index 2)
$ws = $wb.Worksheets.Item(2)
$ws.Range("F9").Value = 103
$ws.Range("F10").Value = 476
$ws.Range("F11").Value = 144579
$ws.Range("F12").Value = 144579
$ws.Range("F18").Value = 228
$ws.Range("F19").Value = 334
$ws.Range("F21").Value = 416
$ws.Range("F22").Value = 416
$ws.Range("F23").Value = 124
$ws.Range("F24").Value = 83
$ws.Range("F27").Value = 553
$ws.Range("F28").Value = 88
$ws.Range("F31").Value = 59
$ws.Range("F32").Value = 344
$ws.Range("F33").Value = 271
$ws.Range("F36").Value = 20
$ws.Range("F38").Value = 192

# sheet3 (index 3)
$ws = $wb.Worksheets.Item(3)
$ws.Range("F4").Value = 3126
$ws.Range("F7").Value = 822
$ws.Range("F8").Value = 1174
$ws.Range("F9").Value = 637
$ws.Range("F10").Value = 1590
$ws.Range("F11").Value = 476
$ws.Range("F12").Value = 86
$ws.Range("F13").Value = 1882

# sheet4 (index 4)
$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 822
$ws.Range("F3").Value = 1174
$ws.Range("F4").Value = 637
$ws.Range("F6").Value = 1590
$ws.Range("F7").Value = 476
$ws.Range("F8").Value = 145
$ws.Range("F9").Value = 2042
$ws.Range("G9").Value = "不可售"
$ws.Range("F10").Value = 86
$ws.Range("F11").Value = 1882
$ws.Range("F12").Value = 4173
$ws.Range("F13").Value = 542
$ws.Range("F14").Value = 659
$ws.Range("F15").Value = 373
$ws.Range("F16").Value = 2195
$ws.Range("F17").Value = 398
$ws.Range("F18").Value = 657931
$ws.Range("F19").Value = 103
$ws.Range("F20").Value = 476
$ws.Range("F21").Value = 1624
$ws.Range("F22").Value = 144579
$ws.Range("F23").Value = 1452
$ws.Range("F24").Value = 670
$ws.Range("F25").Value = 541
$ws.Range("F26").Value = 1269
$ws.Range("F27").Value = 2222
$ws.Range("F28").Value = 1134
$ws.Range("F29").Value = 2691
$ws.Range("F30").Value = 1547
$ws.Range("F31").Value = 801
$ws.Range("F33").Value = 1532
$ws.Range("F34").Value = 416
$ws.Range("F35").Value = 527
$ws.Range("F36").Value = 124
$ws.Range("F37").Value = 1083
$ws.Range("F38").Value = 1078
$ws.Range("F40").Value = 2012
$ws.Range("F41").Value = 1367
$ws.Range("F42").Value = 1222
$ws.Range("F43").Value = 2439
$ws.Range("F44").Value = 1142
$ws.Range("F45").Value = 344
$ws.Range("F46").Value = 344
$ws.Range("F47").Value = 271
$ws.Range("F48").Value = 2575
$ws.Range("F49").Value = 207
$ws.Range("F50").Value = 981
$ws.Range("F51").Value = 3114
$ws.Range("F52").Value = 150
